# OS_test_sorter.xlsx - rebuild the OS1..OS4 blocks shifted from B:C into A:B
# (one row up, one column left) and append three new blocks OS5..OS7 using
# the same layout (header row + D0N/D0P/D1N/D1P rows + blank separator row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old B:C data before laying out the new A:B data.
$ws.Range("B2:C24").ClearContents()

$labels = @("OS1", "OS2", "OS3", "OS4", "OS5", "OS6", "OS7")
$subLabels = @("D0N", "D0P", "D1N", "D1P")

$row = 1
foreach ($label in $labels) {
    $ws.Cells.Item($row, 1).Value = $label
    $row = $row + 1

    foreach ($sub in $subLabels) {
        $ws.Cells.Item($row, 1).Value = $sub
        $ws.Cells.Item($row, 2).Value = 0.234
        $row = $row + 1
    }

    # blank separator row between blocks
    $row = $row + 1
}
